# Updates the cryptos list: refreshed prices and 1h volume percentages
# for each coin row, and swaps the WrappedBTC/Avalanche row order plus
# replaces the MultiversX row with FTXToken.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Price cells as Text so numeric-looking values (e.g. "5.40") keep
# their exact formatting instead of being auto-converted to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the updated values
$ws.Range("D2").Value = "37.162.93"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.054.17"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "250.03"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "0.666"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").Value = "59.82"
$ws.Range("E7").Value = "  +7.69%  "
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "16.18"
$ws.Range("E12").Value = "  +6.38%  "
$ws.Range("D13").Value = "2.352.85"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "0.819"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "5.62"
$ws.Range("E15").Value = "  +7.14%  "
$ws.Range("D16").Value = "2.055.59"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "18.10"
$ws.Range("E17").Value = "  +27.36%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "37.153.17"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "74.94"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "5.40"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "238.34"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  +11.13%  "
$ws.Range("D26").Value = "168.75"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "9.39"
$ws.Range("E27").Value = "  +3.63%  "
$ws.Range("D28").Value = "19.94"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("E30").Value = "  +8.50%  "
$ws.Range("D31").Value = "4.80"
$ws.Range("E31").Value = "  +5.55%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  +4.40%  "
$ws.Range("D34").Value = "0.0897"
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").Value = "2.26"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").Value = "  +8.33%  "
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "5.29"
$ws.Range("E40").Value = "  +31.32%  "
$ws.Range("E41").Value = "  +12.90%  "
$ws.Range("D42").Value = "17.72"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").Value = "96.57"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").Value = "1.288.69"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "6.81"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "2.249.62"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "3.45"
$ws.Range("E51").Value = "  -21.54%  "
